# Update roll numbers from "23BME.." to "22BME.." (new intake) across the
# four award-list sheets (Mid Term, Assignment & Sessional, Practical,
# Final Term). Names/columns B are untouched - only column C (Roll Number)
# changes. Formula-driven columns (L, P, T, X, AB, AF, O, S, W, AA, AE, AI,
# AM, ...) reference $C.. directly and recalc automatically.

$wb = $excel.ActiveWorkbook

$rollNumbers = @(
    "22BME01","22BME02","22BME03","22BME04","22BME05","22BME06","22BME07",
    "22BME08","22BME09","22BME10","22BME11","22BME12","22BME14","22BME15",
    "22BME16","22BME17","22BME18","22BME19","22BME21","22BME22","22BME24",
    "22BME25","22BME26","22BME27","22BME29","22BME31","22BME32","22BME33"
)

function Set-RollNumbers {
    param($SheetName, $StartRow)

    $ws = $wb.Worksheets.Item($SheetName)
    $endRow = $StartRow + $rollNumbers.Length - 1

    # Select the destination range first, mirroring how the edit was made
    # in the live workbook (select C<start>:C<end>, then fill in values) -
    # this is also what drives the saved selection/sqref for the sheet.
    $ws.Range("C" + $StartRow + ":C" + $endRow).Select()

    for ($i = 0; $i -lt $rollNumbers.Length; $i++) {
        $ws.Cells.Item($StartRow + $i, 3).Value = $rollNumbers[$i]
    }

    # Normalize the formatting across the whole filled range to match the
    # anchor cell's style (some rows further down still carried an older
    # alternating style that the bulk fill overwrites uniformly).
    $ws.Range("C" + $StartRow).Copy() | Out-Null
    $ws.Range("C" + $StartRow + ":C" + $endRow).PasteSpecial(-4122) | Out-Null
}

# Mid Term Award: rows 16-43
Set-RollNumbers "Mid Term Award" 16

# Assignment & Sessional: rows 15-42
Set-RollNumbers "Assignment & Sessional" 15

# Practical Award: rows 16-43
Set-RollNumbers "Practical Award" 16

# Final Term Award: rows 16-43
Set-RollNumbers "Final Term Award" 16

# Final Term Award ends up the active sheet/tab, with the cursor left on F19.
$finalSheet = $wb.Worksheets.Item("Final Term Award")
$finalSheet.Activate()
$finalSheet.Range("F19").Select()
